$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels in row 3 to "log(t)"
$ws.Range("H3").Value = "log(t)"
$ws.Range("I3").Value = "log(t)"
$ws.Range("J3").Value = "log(t)"

# Update the active selection to J3
$ws.Range("J3").Select()
